# Update "countries & provincias Spain" data sheet.
# 1) Bolivia moves up in the sorted (by Casos totales desc) list: it now sits
#    right above Ucrania (between Kuwait and Ucrania) with freshly refreshed
#    stats; Ucrania and Emiratos Arabes Unidos each shift down one row and
#    keep the stats that used to belong to the row above them.
# 2) A handful of other country rows get refreshed case counts.
# 3) The "Datos actualizados ..." timestamp string is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 38: Bolivia (new data, moved ahead of Ucrania) ---
$ws.Cells.Item(38, 1).Value = "Bolivia"
$ws.Cells.Item(38, 2).Value = 58138
$ws.Cells.Item(38, 3).Value = 2036
$ws.Cells.Item(38, 4).Value = 18200
$ws.Cells.Item(38, 5).Value = 37832
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 57
$ws.Cells.Item(38, 8).Value = 2106

# --- Row 39: Ucrania (shifted down, keeps its previous stats) ---
$ws.Cells.Item(39, 1).Value = "Ucrania"
$ws.Cells.Item(39, 2).Value = 58111
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = 30525
$ws.Cells.Item(39, 5).Value = 26109
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(39, 8).Value = 1477

# --- Row 40: Emiratos Arabes Unidos (shifted down, keeps its previous stats) ---
$ws.Cells.Item(40, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(40, 2).Value = 56711
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 48917
$ws.Cells.Item(40, 5).Value = 7456
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 338

# --- Row 55: refreshed stats ---
$ws.Cells.Item(55, 2).Value = 32793
$ws.Cells.Item(55, 3).Value = 1048
$ws.Cells.Item(55, 4).Value = 3661
$ws.Cells.Item(55, 5).Value = 28241
$ws.Cells.Item(55, 7).Value = 34
$ws.Cells.Item(55, 8).Value = 891

# --- Row 74: refreshed stats ---
$ws.Cells.Item(74, 2).Value = 11802
$ws.Cells.Item(74, 3).Value = 361
$ws.Cells.Item(74, 4).Value = 8273
$ws.Cells.Item(74, 5).Value = 3407

# --- Row 76: refreshed stats ---
$ws.Cells.Item(76, 2).Value = 11483
$ws.Cells.Item(76, 4).Value = 3972
$ws.Cells.Item(76, 5).Value = 7401
$ws.Cells.Item(76, 8).Value = 110

# --- Row 197: refreshed stats ---
$ws.Cells.Item(197, 2).Value = 38
$ws.Cells.Item(197, 5).Value = 9

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 19 de Julio de 2020 a las 05:25"
